# Corrected name of chapter 3.
# "Integrated experimental strategy" -> "Integrated Experimental Strategy"
# (capitalise the "e" in "experimental" and the "s" in "strategy")

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)            # "TextBox 3" - holds the chapter-3 label
$tr = $sh.TextFrame.TextRange

# Sanity check - make sure we are editing the expected text.
if ($tr.Text -eq "Integrated experimental strategy") {

    $origHeight = $sh.Height

    # 1) Capitalise the "s" of "strategy" (char 25 of the original string).
    $tail = $tr.Characters(25, 1)
    $tail.Text = "S"

    # 2) Re-type the leading portion of the text ("Integrated experimental ")
    #    as "Integrated Experimental " - this both fixes the capital "E" and
    #    merges that whole span back into a single run.
    $head = $tr.Characters(1, 24)
    $head.Text = "Integrated Experimental "

    # The textbox auto-fits its height to the text; restore the original
    # box height so only the wording changes.
    $sh.Height = $origHeight
}

Write-Host $tr.Text
